$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 118517.82
$ws.Range("I28").Value = 154485
$ws.Range("J28").Value = 1624.5
$ws.Range("K28").Value = 154485
$ws.Range("L28").Value = 1624.5
$ws.Range("M28").Value = -154000
$ws.Range("N28").Value = -2594.5
$ws.Range("H42").Value = 162.42857
$ws.Range("J42").Value = 232.5
$ws.Range("L42").Value = 697.5
$ws.Range("N42").Value = -1157.5
$ws.Range("H43").Value = 6999.8
$ws.Range("J43").Value = 8333.333000000001
$ws.Range("L43").Value = 8333.333000000001
$ws.Range("N43").Value = -8471.333000000001
$ws.Range("H69").Value = 10015
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H70").Value = 35716384
$ws.Range("I70").Value = 1186
$ws.Range("J70").Value = 41668916
$ws.Range("K70").Value = 3558
$ws.Range("L70").Value = 125006748
$ws.Range("M70").Value = -3288
$ws.Range("N70").Value = -125007288
$ws.Range("H72").Value = 10015
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H73").Value = 35716384
$ws.Range("I73").Value = 1186
$ws.Range("J73").Value = 41668916
$ws.Range("K73").Value = 3558
$ws.Range("L73").Value = 125006748
$ws.Range("M73").Value = -2622
$ws.Range("N73").Value = -125008620
$ws.Range("H86").Value = 86585944
$ws.Range("I86").Value = 62510660
$ws.Range("K86").Value = 62510660
$ws.Range("M86").Value = -62509537
$ws.Range("H89").Value = 86585944
$ws.Range("I89").Value = 62510660
$ws.Range("K89").Value = 312553300
$ws.Range("M89").Value = -312547684
$ws.Range("H92").Value = 15625767
$ws.Range("I92").Value = 19231572
$ws.Range("K92").Value = 19231572
$ws.Range("M92").Value = -19230324
$ws.Range("H106").Value = 9527211
$ws.Range("I106").Value = 9527211
$ws.Range("K106").Value = 9527211
$ws.Range("M106").Value = -9526580
$ws.Range("H137").Value = 36220.12
$ws.Range("I137").Value = 54319.375
$ws.Range("J137").Value = 4043.6667
$ws.Range("K137").Value = 162958.125
$ws.Range("L137").Value = 12131.0001
$ws.Range("M137").Value = -160408.125
$ws.Range("N137").Value = -17231.0001
$ws.Range("H138").Value = 1901.9899
$ws.Range("J138").Value = 2472.0806
$ws.Range("L138").Value = 7416.2418
$ws.Range("N138").Value = -17696.2418

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3953.4092
$ws.Range("I63").Value = 2330.4443
$ws.Range("K63").Value = 2330.4443
$ws.Range("M63").Value = -1644.4443
$ws.Range("H66").Value = 3953.4092
$ws.Range("I66").Value = 2330.4443
$ws.Range("K66").Value = 11652.2215
$ws.Range("M66").Value = -8220.2215
$ws.Range("H74").Value = 2538.0513
$ws.Range("I74").Value = 2389.3057
$ws.Range("J74").Value = 4323
$ws.Range("K74").Value = 2389.3057
$ws.Range("L74").Value = 4323
$ws.Range("M74").Value = -1515.3057
$ws.Range("N74").Value = -6071
$ws.Range("H77").Value = 2538.0513
$ws.Range("I77").Value = 2389.3057
$ws.Range("J77").Value = 4323
$ws.Range("K77").Value = 11946.5285
$ws.Range("L77").Value = 21615
$ws.Range("M77").Value = -7578.5285
$ws.Range("N77").Value = -30351
$ws.Range("H132").Value = 246425.39
$ws.Range("I132").Value = 296323.06
$ws.Range("J132").Value = 4065.2856
$ws.Range("K132").Value = 888969.1799999999
$ws.Range("L132").Value = 12195.8568
$ws.Range("M132").Value = -886439.1799999999
$ws.Range("N132").Value = -17255.8568

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5120.242
$ws.Range("I31").Value = 2174.4856
$ws.Range("J31").Value = 8446.097
$ws.Range("K31").Value = 2174.4856
$ws.Range("L31").Value = 8446.097
$ws.Range("M31").Value = -1879.4856
$ws.Range("N31").Value = -9036.097
$ws.Range("H34").Value = 5120.242
$ws.Range("I34").Value = 2174.4856
$ws.Range("J34").Value = 8446.097
$ws.Range("K34").Value = 2174.4856
$ws.Range("L34").Value = 8446.097
$ws.Range("M34").Value = -1972.4856
$ws.Range("N34").Value = -8850.097
$ws.Range("H55").Value = 9980
$ws.Range("I55").Value = 9980
$ws.Range("K55").Value = 9980
$ws.Range("M55").Value = -9665
$ws.Range("H58").Value = 2795.625
$ws.Range("I58").Value = 2679.6052
$ws.Range("K58").Value = 2679.6052
$ws.Range("M58").Value = -2476.6052
$ws.Range("H132").Value = 4349.5
$ws.Range("I132").Value = 4339.773
$ws.Range("J132").Value = 4376.25
$ws.Range("K132").Value = 13019.319
$ws.Range("L132").Value = 13128.75
$ws.Range("M132").Value = -10489.319
$ws.Range("N132").Value = -18188.75
$ws.Range("H136").Value = 2795.625
$ws.Range("I136").Value = 2679.6052
$ws.Range("K136").Value = 8038.8156
$ws.Range("M136").Value = -5488.8156

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9421.714
$ws.Range("I113").Value = 2396.375
$ws.Range("K113").Value = 2396.375
$ws.Range("M113").Value = -226.375
$ws.Range("H120").Value = 79891.664
$ws.Range("J120").Value = 79891.664
$ws.Range("L120").Value = 79891.664
$ws.Range("N120").Value = -89567.664
$ws.Range("H122").Value = 1078
$ws.Range("I122").Value = 1117.5
$ws.Range("K122").Value = 3352.5
$ws.Range("M122").Value = -902.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1353.8
$ws.Range("I61").Value = 948.6667
$ws.Range("K61").Value = 948.6667
$ws.Range("M61").Value = -746.6667
$ws.Range("H68").Value = 9654.362999999999
$ws.Range("I68").Value = 10939
$ws.Range("J68").Value = 8583.833000000001
$ws.Range("K68").Value = 10939
$ws.Range("L68").Value = 8583.833000000001
$ws.Range("M68").Value = -10190
$ws.Range("N68").Value = -10081.833
$ws.Range("H71").Value = 9654.362999999999
$ws.Range("I71").Value = 10939
$ws.Range("J71").Value = 8583.833000000001
$ws.Range("K71").Value = 54695
$ws.Range("L71").Value = 42919.165
$ws.Range("M71").Value = -50951
$ws.Range("N71").Value = -50407.165
$ws.Range("H81").Value = 95989
$ws.Range("J81").Value = 95989
$ws.Range("L81").Value = 95989
$ws.Range("N81").Value = -97985
$ws.Range("H84").Value = 95989
$ws.Range("J84").Value = 95989
$ws.Range("L84").Value = 287967
$ws.Range("N84").Value = -297951
$ws.Range("H113").Value = 1353.8
$ws.Range("I113").Value = 948.6667
$ws.Range("K113").Value = 948.6667
$ws.Range("M113").Value = 1221.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5363.2856
$ws.Range("H65").Value = 5363.2856
$ws.Range("H81").Value = 75925.86
$ws.Range("I81").Value = 146223.42
$ws.Range("K81").Value = 292446.84
$ws.Range("M81").Value = -291385.84
$ws.Range("H84").Value = 75925.86
$ws.Range("I84").Value = 146223.42
$ws.Range("K84").Value = 1462234.2
$ws.Range("M84").Value = -1456930.2
$ws.Range("H126").Value = 5019.3125
$ws.Range("I126").Value = 5801
$ws.Range("J126").Value = 3299.6
$ws.Range("K126").Value = 17403
$ws.Range("L126").Value = 9898.799999999999
$ws.Range("M126").Value = -14933
$ws.Range("N126").Value = -14838.8
